$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B21").Value = "5.2 Perform Statement Coverage Testing"
$ws.Range("B21").Characters(4,35).Font.Bold = $true
$ws.Range("B21").Characters(4,35).Font.Color = 255
$ws.Range("B21").Characters(4,35).Font.Name = "Calibri"
$ws.Range("B21").Characters(4,35).Font.Size = 12

$ws.Range("D19").Font.Color = 255
$ws.Range("D19").Font.Name = "Calibri"
$ws.Range("D19").Font.Size = 12
